$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data describing the "CB construction" benchmark entry
$ws.Range("C26").Value = "CB construction (2000 sample, 16 CW)"
$ws.Range("D26").Value = 456541351
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Formula = "=D26/225000000*1000"
$ws.Range("F26").Formula = "=E26/10"
$ws.Range("F26").NumberFormat = "0.00%"

# Update the selected cell to reflect where the author ended up editing
[void]$ws.Range("C31").Select()
